$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.952.78"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "1.893.54"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +1.66%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "335.70"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("E6").Value = "  +1.40%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4698"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.58%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3935"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.26%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "47.49"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.43%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.08072"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.05%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.024"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.28%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "21.88"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "1.862.88"
$ws.Range("E13").Value = "  -0.63%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.983"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.28%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.151"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.69%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.016"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.52%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.06793"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +3.44%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "87.43"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.00001047"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.76%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "17.25"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.11%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "1.014"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("D22").Value = "27.990.32"
$ws.Range("E22").Value = "  +1.11%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.532"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.27%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.03"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.09%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.347"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").Value = "2.104.98"
$ws.Range("E26").Value = "  +0.36%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "159.81"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.43%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "20.09"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.97%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.097"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.46%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "5.492"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.92%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "122.11"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.50%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.9773"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.17%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.09516"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.642"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.81%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.410"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -4.66%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.391"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.39%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.06145"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.44%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02265"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.54%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.225"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.05%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "8.090"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.96%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.6017"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.17%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1892"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.39%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "10.31"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("E44").Value = "  +1.55%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.5716"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "12.20"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.42%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.406"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.24%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.943"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.02%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.06936"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "114.11"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.072"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.58%  "
